# Generate Report for Handoff
# Updates the localization-status workbook: file "26763ceb-...md" is renamed /
# replaced by "6a1f3617-...md" (handed off, not yet handed back) and file
# "e66f312d-...md" is replaced by "ffff09a99ca6-...md" (a content-duplicate of
# the first file, so it shares its handoff artifact).

$wb = $excel.ActiveWorkbook

# ---- new identifiers -------------------------------------------------
$file1Md      = "6a1f3617-07d8-4b4c-add0-535573617591.md"
$file1Path    = "e2e\6a1f3617-07d8-4b4c-add0-535573617591.md"
$file2Md      = "ffff09a99ca6-b8cb-4b53-83b7-9f1935683d43.md"
$file2Path    = "e2e\ffff09a99ca6-b8cb-4b53-83b7-9f1935683d43.md"

$status       = "Ready for handoff"
$hoDate       = "2016-08-29 07:04:05"

$zhcnXlf      = "6a1f3617-07d8-4b4c-add0-535573617591.21614cdb91720701de0327b3fbcce43157ad5500.zh-cn.xlf"
$dedeXlf      = "6a1f3617-07d8-4b4c-add0-535573617591.21614cdb91720701de0327b3fbcce43157ad5500.de-de.xlf"
$handoffDtZh  = "2016-08-29 07:03:57"
$handoffDtDe  = $hoDate
$handbackDt0  = "0001-01-01 00:00:00"

# =======================================================================
# Sheet "Overview"
# =======================================================================
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = $file1Md
$ov.Range("A3").Value = $file2Md

foreach ($h in $ov.Hyperlinks) {
    if ($h.Range.Address() -eq '$B$2') { $h.TextToDisplay = $file1Path }
    if ($h.Range.Address() -eq '$B$3') { $h.TextToDisplay = $file2Path }
}
$ov.Range("B2").Value = $file1Path
$ov.Range("B3").Value = $file2Path

$ov.Range("E2").Value = $status
$ov.Range("F2").Value = $status
$ov.Range("E3").Value = $status
$ov.Range("F3").Value = $status

$ov.Range("G2").Value = $hoDate
$ov.Range("G3").Value = $hoDate

$ov.Columns.Item(5).AutoFit() | Out-Null
$ov.Columns.Item(6).AutoFit() | Out-Null

# =======================================================================
# Sheet "zh-cn"
# =======================================================================
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = $file1Md
$zh.Range("A3").Value = $file2Md
foreach ($h in $zh.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') { $h.TextToDisplay = $file1Md }
    if ($h.Range.Address() -eq '$A$3') { $h.TextToDisplay = $file2Md }
}

$zh.Range("C2").Value = $status
$zh.Range("C3").Value = $status

# File 2 is now flagged as a content duplicate of File 1.
$zh.Range("F3").Value = "True"

# Both rows now reference the (single) handoff package generated for File 1.
$zh.Range("G2").Value = $zhcnXlf
$zh.Range("G3").Value = $zhcnXlf

$zh.Range("H2").Value = $handoffDtZh
$zh.Range("H3").Value = $handoffDtZh

# Not handed back yet -> target/handback file columns are blank and the
# handback hyperlink is removed.
foreach ($h in @($zh.Hyperlinks)) {
    if ($h.Range.Address() -eq '$I$2' -or $h.Range.Address() -eq '$I$3') { $h.Delete() }
}
$zh.Range("I2").Value = ""
$zh.Range("I2").Style = "Normal"
$zh.Range("I3").Value = ""
$zh.Range("I3").Style = "Normal"

$zh.Range("J2").Value = ""
$zh.Range("J3").Value = ""

$zh.Range("K2").Value = $handbackDt0
$zh.Range("K3").Value = $handbackDt0

$zh.Columns.Item(3).AutoFit() | Out-Null
$zh.Columns.Item(9).AutoFit() | Out-Null
$zh.Columns.Item(10).AutoFit() | Out-Null

# =======================================================================
# Sheet "de-de"
# =======================================================================
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = $file1Md
$de.Range("A3").Value = $file2Md
foreach ($h in $de.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') { $h.TextToDisplay = $file1Md }
    if ($h.Range.Address() -eq '$A$3') { $h.TextToDisplay = $file2Md }
}

$de.Range("C2").Value = $status
$de.Range("C3").Value = $status

$de.Range("F3").Value = "True"

$de.Range("G2").Value = $dedeXlf
$de.Range("G3").Value = $dedeXlf

$de.Range("H2").Value = $handoffDtDe
$de.Range("H3").Value = $handoffDtDe

foreach ($h in @($de.Hyperlinks)) {
    if ($h.Range.Address() -eq '$I$2' -or $h.Range.Address() -eq '$I$3') { $h.Delete() }
}
$de.Range("I2").Value = ""
$de.Range("I2").Style = "Normal"
$de.Range("I3").Value = ""
$de.Range("I3").Style = "Normal"

$de.Range("J2").Value = ""
$de.Range("J3").Value = ""

$de.Range("K2").Value = $handbackDt0
$de.Range("K3").Value = $handbackDt0

$de.Columns.Item(3).AutoFit() | Out-Null
$de.Columns.Item(9).AutoFit() | Out-Null
$de.Columns.Item(10).AutoFit() | Out-Null
